$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new column at AM (pushes the existing "assign" column from AM to AN)
$ws.Columns("AM").Insert()

# 2. New header cell AM1: "fit"
$hdr = $ws.Range("AM1")
$hdr.Value2 = "fit"
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.Borders.Item(8).LineStyle = -4142   # xlEdgeTop -> none
$hdr.Borders.Item(9).LineStyle = -4142   # xlEdgeBottom -> none
$left = $hdr.Borders.Item(7)             # xlEdgeLeft
$left.LineStyle = 1                      # xlContinuous
$left.Weight = 2                         # xlThin
$right = $hdr.Borders.Item(10)           # xlEdgeRight
$right.LineStyle = 1
$right.Weight = 2

# 3. Fill in the new "fit" cluster values for rows 2-51
$fitValues = @(2,2,2,2,2,2,2,2,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,4,1,1,1,0,0,0,0,0,0,0,0,0,0,3)
for ($i = 0; $i -lt $fitValues.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 39).Value2 = $fitValues[$i]
}

# 4. Match the new column's width to its left neighbor so the <cols> entries coalesce
$ws.Columns("AL").ColumnWidth = $ws.Columns("AL").ColumnWidth
$ws.Columns("AM").ColumnWidth = $ws.Columns("AL").ColumnWidth

# 5. Re-point the AutoFilter over the full, now-wider range
$ws.AutoFilterMode = $false
$ws.Range("A1:AN51").AutoFilter()

# 6. Keep the hidden _FilterDatabase defined name in sync with the AutoFilter range
for ($i = 1; $i -le $wb.Names.Count; $i++) {
    $n = $wb.Names.Item($i)
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$AN`$51"
    }
}

# 7. Restore the selection to match the edited workbook
$ws.Range("AK3").Select()

Write-Host "done"
